# ---------------------------------------------------------------------------
# Restructure NinzaData.xlsx:
#   Sheet1 / NinzaAutomation  ->  Campaigns, Contacts, Products
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ============================================================
# Phase 1 - build the sheet skeleton.
# The order sheets are added/removed in controls both the internal
# sheetId counter (next = current-max-id + 1) and final tab order, so it
# is deliberately sequenced:
#   1) rename the existing Sheet1 -> Contacts (keeps its original id)
#   2) add Campaigns while NinzaAutomation (id 2) still exists -> id 3
#   3) drop NinzaAutomation
#   4) add Products now that the highest live id is 3 -> id 4
#   5) slide Products to the end, after Contacts
# ============================================================
$wsContacts = $wb.Worksheets.Item(1)
$wsContacts.Name = "Contacts"

$wsCampaigns = $wb.Worksheets.Add()
$wsCampaigns.Name = "Campaigns"

$wsNinza = $wb.Worksheets.Item("NinzaAutomation")
$wsNinza.Delete() | Out-Null

$wsProducts = $wb.Worksheets.Add()
$wsProducts.Name = "Products"
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# ============================================================
# Phase 2 - fill in the cell content.
# Sheet objects are re-fetched by name after every structural change
# above (Add/Delete/Move can rebind stale worksheet variables), and the
# write order below is what drives the shared-string table order.
# ============================================================

# --- Contacts header row ---
$wsContacts = $wb.Worksheets.Item("Contacts")
$wsContacts.Range("A1").Value = "Organization"
$wsContacts.Range("B1").Value = "Title"
$wsContacts.Range("C1").Value = "Contact Name"
$wsContacts.Range("D1").Value = "Mobile"

# --- Products header row (Vendor added later, below) ---
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("A1").Value = "ProductName"
$wsProducts.Range("B1").Value = "Category"
$wsProducts.Range("C1").Value = "Quantity"
$wsProducts.Range("D1").Value = "PricePerUnit"

# --- Campaigns header row ---
$wsCampaigns = $wb.Worksheets.Item("Campaigns")
$wsCampaigns.Range("A1").Value = "CampaignName"
$wsCampaigns.Range("B1").Value = "TargetSize"

# --- back to Products: extra column ---
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("E1").Value = "Vendor"

# --- back to Campaigns: sample data row (TargetSize kept as text, like a
#     quote-prefixed "10" typed in Excel) ---
$wsCampaigns = $wb.Worksheets.Item("Campaigns")
$wsCampaigns.Range("A2").Value = "Qspiders-4510"
$wsCampaigns.Range("B2").Value = "'10"

# ============================================================
# Phase 3 - column widths (best effort - manually sized like a user
# dragging/typing a column width) and per-sheet cursor position.
# ============================================================
$wsCampaigns = $wb.Worksheets.Item("Campaigns")
$wsCampaigns.Columns.Item(1).ColumnWidth = 13.166666666666666
$wsCampaigns.Columns.Item(2).ColumnWidth = 8.5
$wsCampaigns.Range("B3").Select() | Out-Null

$wsContacts = $wb.Worksheets.Item("Contacts")
$wsContacts.Columns.Item(1).ColumnWidth = 12
$wsContacts.Columns.Item(3).ColumnWidth = 14.5
$wsContacts.Range("B6").Select() | Out-Null

$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Columns.Item(1).ColumnWidth = 11.333333333333334
$wsProducts.Columns.Item(4).ColumnWidth = 10.833333333333334
$wsProducts.Range("E1").Select() | Out-Null

# ============================================================
# Phase 4 - Campaigns is the first/active tab, as in the source file.
# ============================================================
$wsCampaigns = $wb.Worksheets.Item("Campaigns")
$wsCampaigns.Select() | Out-Null
